$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 299, shifting existing rows 299:321 down to 300:322.
$ws.Rows.Item(299).Insert()

# Populate the new row 299 with the weekly record for "Femacal de La Calera - Ciboulette".
$ws.Cells.Item(299, 1).Value = 3
$ws.Cells.Item(299, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(299, 3).Value = "Coquimbo"
$ws.Cells.Item(299, 4).Value = 44714
$ws.Cells.Item(299, 5).Value = 5
$ws.Cells.Item(299, 6).Value = 100112039
$ws.Cells.Item(299, 7).Value = "Ciboulette"
$ws.Cells.Item(299, 8).Value = "Sin especificar"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 120
$ws.Cells.Item(299, 11).Value = 1500
$ws.Cells.Item(299, 12).Value = 1500
$ws.Cells.Item(299, 13).Value = 1500
$ws.Cells.Item(299, 14).Value = "$/docena de atados"
$ws.Cells.Item(299, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(299, 16).Value = 500
$ws.Cells.Item(299, 17).Value = 3
$ws.Cells.Item(299, 18).Value = "Hortaliza"
